$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43 (Ordem 42 / Codigo 2720) : EEB PREF AVELINO MULLER ---
# Copy formatting from an existing "school name" cell (style s="1") before
# writing the new text so no new style entries are introduced.
$ws.Range("C42").Copy($ws.Range("C43"))
$ws.Range("C43").Value = "EEB PREF AVELINO MULLER"

$ws.Range("D43").Value = "x"
$ws.Range("E43").Value = "x"
$ws.Range("H43").Value = 20
$ws.Range("I43").Value = 47
$ws.Range("J43").Value = 23
$ws.Range("K43").Value = 50
$ws.Range("L43").Value = 28
$ws.Range("M43").Value = 3000
$ws.Range("N43").Value = 50
$ws.Range("O43").Value = 555
$ws.Range("P43").Value = 3
$ws.Range("Q43").Value = "OK"

# --- Row 44 (Ordem 43 / Codigo 2674) : EEB DR ADERBAL RAMOS DA SILVA ---
$ws.Range("B44").Value = 2674

$ws.Range("C42").Copy($ws.Range("C44"))
$ws.Range("C44").Value = "EEB DR ADERBAL RAMOS DA SILVA"

$ws.Range("D44").Value = "x"
$ws.Range("E44").Value = "x"
$ws.Range("F44").Value = "x"
$ws.Range("G44").Value = "x"
$ws.Range("H44").Value = 26
$ws.Range("I44").Value = 47
$ws.Range("J44").Value = 22
$ws.Range("K44").Value = 50
$ws.Range("L44").Value = 15
$ws.Range("M44").Value = 3550
$ws.Range("N44").Value = 592
$ws.Range("O44").Value = 585
$ws.Range("P44").Value = 4
$ws.Range("Q44").Value = "OK"

# --- Row 45 (Ordem 44) : Codigo 2739 ---
$ws.Range("B45").Value = 2739

# --- Update the active sheet view / selection ---
$excel.ActiveWindow.ScrollRow = 31
[void]$ws.Range("C45").Select()
